$wb = $excel.ActiveWorkbook

# --- Update "Training Data" sheet: column D (traffic_volume) values ---
$wsTrain = $wb.Worksheets.Item("Training Data")

$wsTrain.Range("D2").Value = 16202
$wsTrain.Range("D3").Value = 11801
$wsTrain.Range("D4").Value = 11055
$wsTrain.Range("D5").Value = 39399
$wsTrain.Range("D6").Value = 32392
$wsTrain.Range("D7").Value = 33049
$wsTrain.Range("D8").Value = 42524
$wsTrain.Range("D9").Value = 36827
$wsTrain.Range("D10").Value = 40086
$wsTrain.Range("D11").Value = 41062
$wsTrain.Range("D12").Value = 40610
$wsTrain.Range("D13").Value = 50361
$wsTrain.Range("D14").Value = 54110
$wsTrain.Range("D15").Value = 34279
$wsTrain.Range("D16").Value = 42675
$wsTrain.Range("D17").Value = 40263
$wsTrain.Range("D18").Value = 36358
$wsTrain.Range("D19").Value = 38424
$wsTrain.Range("D20").Value = 44354
$wsTrain.Range("D21").Value = 38732
$wsTrain.Range("D22").Value = 39234
$wsTrain.Range("D23").Value = 42488
$wsTrain.Range("D24").Value = 41926
$wsTrain.Range("D25").Value = 50080
$wsTrain.Range("D26").Value = 52926
$wsTrain.Range("D27").Value = 40132
$wsTrain.Range("D28").Value = 28390
$wsTrain.Range("D29").Value = 7710
$wsTrain.Range("D30").Value = 19584
$wsTrain.Range("D31").Value = 24174
$wsTrain.Range("D32").Value = 29230
$wsTrain.Range("D33").Value = 32758
$wsTrain.Range("D34").Value = 37779
$wsTrain.Range("D35").Value = 44243
$wsTrain.Range("D36").Value = 37700
$wsTrain.Range("D37").Value = 42238
$wsTrain.Range("D38").Value = 42036
$wsTrain.Range("D39").Value = 34395
$wsTrain.Range("D40").Value = 8996
$wsTrain.Range("D41").Value = 28242
$wsTrain.Range("D42").Value = 34245
$wsTrain.Range("D43").Value = 35938
$wsTrain.Range("D44").Value = 42164
$wsTrain.Range("D45").Value = 37944
$wsTrain.Range("D46").Value = 40373
$wsTrain.Range("D47").Value = 43586
$wsTrain.Range("D48").Value = 39308
$wsTrain.Range("D49").Value = 48682
$wsTrain.Range("D50").Value = 46705
$wsTrain.Range("D51").Value = 34679
$wsTrain.Range("D52").Value = 41350
$wsTrain.Range("D53").Value = 44583
$wsTrain.Range("D54").Value = 37851
$wsTrain.Range("D55").Value = 37327
$wsTrain.Range("D56").Value = 45732
$wsTrain.Range("D57").Value = 39369
$wsTrain.Range("D58").Value = 41948
$wsTrain.Range("D59").Value = 44723
$wsTrain.Range("D60").Value = 39863
$wsTrain.Range("D61").Value = 50245

# --- Clear "Testing Data" sheet: column D (traffic_volume) values for rows 2-13 ---
$wsTest = $wb.Worksheets.Item("Testing Data")
$wsTest.Range("D2:D13").ClearContents()

